# ---------------------------------------------------------------------------
# Applies the "fixed set expire error, add redpacket message" edit to the
# single-slide web3chat deck:
#   1. Shrinks/repositions the round-rect badge background (shape id 4).
#   2. Deletes the chat-bubble picture (shape id 6 / "图形 5").
#   3. Shrinks/repositions the "web3chat" label (shape id 7) and reduces its
#      font size from 54pt to 32pt.
#   4. Bumps the cached datetimeFigureOut field text (2022/3/21 -> 2022/3/26)
#      wherever it is cached (slide master + all slide layouts).
# ---------------------------------------------------------------------------

# Shape Left/Top/Width/Height are exposed as 32-bit floats (points), so a
# naive EMU/12700 round-trip can land 1 EMU short of the target after the
# float32 truncation PowerPoint's COM layer performs when it re-derives EMUs
# from the point value. Binary-search the smallest point value that survives
# the float32 cast and still floors back to the exact target EMU.
#
# NOTE: these helpers must always be invoked *positionally* - calling a
# `param()`-block function with `-Name value` style arguments silently drops
# the bound values (and COM object references!) in this interpreter.
function ConvertTo-ExactPoints {
    param($TargetEmu)
    $base = [double]$TargetEmu / 12700.0
    $lo = 0.0
    $hi = 2.0 / 12700.0
    for ($i = 0; $i -lt 60; $i++) {
        $mid = ($lo + $hi) / 2.0
        $candidate = $base + $mid
        $asSingle = [float]$candidate
        $emu = [math]::Floor([double]$asSingle * 12700.0)
        if ($emu -ge $TargetEmu) {
            $hi = $mid
        } else {
            $lo = $mid
        }
    }
    return $base + $hi
}

function Set-ShapeExactPosition {
    param($Shape, $LeftEmu, $TopEmu, $WidthEmu, $HeightEmu)
    $Shape.Left = ConvertTo-ExactPoints $LeftEmu
    $Shape.Top = ConvertTo-ExactPoints $TopEmu
    $Shape.Width = ConvertTo-ExactPoints $WidthEmu
    $Shape.Height = ConvertTo-ExactPoints $HeightEmu
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Round-rect badge background (id 4, "矩形: 圆角 3") ------------------
$roundRect = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Id -eq 4) { $roundRect = $s.Shapes.Item($i) }
}
Set-ShapeExactPosition $roundRect 4301269 2896829 2122518 674605

# --- 2. Delete the chat-bubble picture (id 6, "图形 5") ---------------------
$pic = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Id -eq 6) { $pic = $s.Shapes.Item($i) }
}
if ($pic -ne $null) { $pic.Delete() }

# --- 3. "web3chat" label (id 7, "矩形 6") -----------------------------------
$label = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Id -eq 7) { $label = $s.Shapes.Item($i) }
}
Set-ShapeExactPosition $label 4377470 2941745 1970117 584775
$label.TextFrame.TextRange.Font.Size = 32

# --- 4. Bump the cached datetimeFigureOut text (slide master + layouts) ----
function Update-CachedDateField {
    param($ShapeCollection)
    for ($i = 1; $i -le $ShapeCollection.Count; $i++) {
        $shp = $ShapeCollection.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "2022/3/21") {
                $tr.Text = "2022/3/26"
            }
        }
    }
}

Update-CachedDateField $p.SlideMaster.Shapes
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-CachedDateField $layouts.Item($li).Shapes
}
